# Weekly update: insert a new price-report row for Acelga (Terminal
# Hortofrutícola Agro Chillán) above the existing row 151, pushing the
# former rows 151-161 down to 152-162.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151; this shifts rows 151:161 down
# to 152:162 (and Excel copies the formatting of the row above, matching
# the date-style used in column D).
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the latest weekly record.
$ws.Cells.Item(151, 1).Value = 7
$ws.Cells.Item(151, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(151, 3).Value = "Ñuble"
$ws.Cells.Item(151, 4).Value = 44461
$ws.Cells.Item(151, 5).Value = 16
$ws.Cells.Item(151, 6).Value = 100112009
$ws.Cells.Item(151, 7).Value = "Acelga"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 120
$ws.Cells.Item(151, 11).Value = 350
$ws.Cells.Item(151, 12).Value = 400
$ws.Cells.Item(151, 13).Value = 375
$ws.Cells.Item(151, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(151, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(151, 16).Value = 375
$ws.Cells.Item(151, 17).Value = 1
$ws.Cells.Item(151, 18).Value = "Hortaliza"
